$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the two rich-text header cells (week number & date range) ---
$ws.Range("A8").Value = "Volume 30   Number  51"
$ws.Range("C9").Value = "Report Covering the Week  12/18/2023  Through  12/24/2023"

# --- Helper functions to change a cell while forcing a specific pre-existing style ---
function Set-TextCellWithStyle($ref, $styleDonor, $text) {
    $ws.Range($styleDonor).Copy()
    $ws.Range($ref).PasteSpecial(-4122)
    $ws.Range($ref).NumberFormat = "@"
    $ws.Range($ref).Value = $text
    $ws.Range($styleDonor).Copy()
    $ws.Range($ref).PasteSpecial(-4122)
}

function Set-NumCellWithStyle($ref, $styleDonor, $val) {
    $ws.Range($styleDonor).Copy()
    $ws.Range($ref).PasteSpecial(-4122)
    $ws.Range($ref).Value = $val
}

# --- Cells that flip between a numeric value and the text placeholder "0" (style 14) ---
Set-TextCellWithStyle "C15" "C14" "0"
Set-TextCellWithStyle "C18" "C14" "0"
Set-TextCellWithStyle "C26" "C14" "0"
Set-TextCellWithStyle "C27" "C14" "0"

# --- Cells that flip from the text placeholder to a real number (style 15 / 16) ---
Set-NumCellWithStyle "C22" "F14" 2
Set-NumCellWithStyle "D22" "F14" 1
Set-NumCellWithStyle "C28" "F14" 1
Set-NumCellWithStyle "C29" "F14" 1
Set-NumCellWithStyle "E22" "K14" 100

# --- Remaining numeric-only updates (style/type unchanged) ---
$map = @{
    "N14" = -86.666666666666
    "N15" = -80.909090909090
    "C16" = 6
    "E16" = 500
    "G16" = 15
    "H16" = 6.666666666666
    "I16" = 164
    "J16" = 199
    "K16" = -17.587939698492
    "L16" = -8.379888268156
    "M16" = -58.585858585858
    "N16" = -92.325690219934
    "C17" = 9
    "D17" = 8
    "E17" = 12.5
    "F17" = 27
    "G17" = 25
    "H17" = 8
    "I17" = 335
    "J17" = 350
    "K17" = -4.285714285714
    "L17" = -3.458213256484
    "M17" = -6.944444444444
    "N17" = -60.818713450292
    "E18" = -100
    "F18" = 5
    "G18" = 11
    "H18" = -54.545454545454
    "J18" = 184
    "K18" = -34.782608695652
    "L18" = -41.176470588235
    "M18" = -62.616822429906
    "N18" = -95.802728226652
    "C19" = 11
    "D19" = 6
    "E19" = 83.333333333333
    "F19" = 38
    "G19" = 35
    "H19" = 8.571428571428
    "I19" = 498
    "J19" = 502
    "K19" = -0.796812749003
    "L19" = -1.386138613861
    "M19" = -25.337331334332
    "N19" = -54.972875226039
    "C20" = 1
    "D20" = 2
    "E20" = -50
    "F20" = 16
    "H20" = 14.285714285714
    "I20" = 120
    "J20" = 120
    "L20" = -5.511811023622
    "M20" = -34.782608695652
    "N20" = -94.910941475827
    "C21" = 27
    "E21" = 42.105263157894
    "F21" = 106
    "G21" = 102
    "H21" = 3.921568627450
    "I21" = 1264
    "J21" = 1387
    "K21" = -8.868060562364
    "L21" = -9.064748201438
    "M21" = -35.477284328739
    "N21" = -86.652587117212
    "F22" = 3
    "G22" = 3
    "I22" = 24
    "J22" = 14
    "K22" = 71.428571428571
    "L22" = 33.333333333333
    "M22" = -11.111111111111
    "C24" = 34
    "D24" = 54
    "E24" = -37.037037037037
    "F24" = 131
    "G24" = 197
    "H24" = -33.502538071066
    "I24" = 1798
    "J24" = 1552
    "K24" = 15.850515463917
    "L24" = 30.478955007256
    "M24" = 34.681647940074
    "C25" = 15
    "D25" = 13
    "E25" = 15.384615384615
    "F25" = 46
    "G25" = 43
    "H25" = 6.976744186046
    "I25" = 603
    "J25" = 559
    "K25" = 7.871198568872
    "L25" = 6.537102473498
    "M25" = -21.688311688311
    "L26" = -14.285714285714
    "E27" = -100
    "G27" = 4
    "H27" = -50
    "J27" = 75
    "K27" = -26.666666666666
    "I28" = 16
    "K28" = 14.285714285714
    "L28" = -15.789473684210
    "M28" = -54.285714285714
    "N28" = -88.652482269503
    "I29" = 14
    "K29" = 7.692307692307
    "L29" = -6.666666666666
    "M29" = -53.333333333333
    "N29" = -87.826086956521
    "F30" = 1
    "G30" = 1
    "H30" = 0
}
foreach ($key in $map.Keys) {
    $ws.Range($key).Value = $map[$key]
}
